# Insert a new "合肥·风禾尽起代号鸢ONLY" row (2024-07-14) as row 7 into the
# two worksheets that list exhibitions ("展览" and "全部类型"), pushing the
# existing rows from 7 downward, and apply the handful of "想去人数" (F
# column) refreshes that came along with the same scrape.

$wb = $excel.ActiveWorkbook

$targetSheets = @("展览", "全部类型")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- insert the new row, shifting everything at/after row 7 down ---
    $ws.Rows.Item(7).Insert()

    # Column A carries a bordered/bold style on every data row; Insert()
    # does not propagate it to the freshly-created row, so copy it over
    # from the row that used to be 7 (now sitting at row 8).
    $ws.Range("A8").Copy()
    $ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false

    # Sequential index column (same numbering scheme as every other row).
    $ws.Range("A7").Value = 6

    # --- populate the new row's data ---
    $ws.Range("B7").Value = "2024-07-14"
    $ws.Range("C7").Value = "合肥·风禾尽起代号鸢ONLY"
    $ws.Range("D7").Value = "长江180艺术街区1-2号楼 圣拉维一站式婚礼宴会艺术中心"
    $ws.Range("E7").Value = "2024.07.14 10:00-07.14 18:00"
    $ws.Range("F7").Value = 0
    $ws.Range("G7").Value = 78
    $ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=88574"
    $ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202407/w6CXgiGo1719918717900.jpeg"

    # --- refresh "想去人数" (F column) counts for events whose numbers moved ---
    # 合肥·星光次元动漫文化节 (row 5, unaffected by the insert)
    $ws.Range("F5").Value = 6741

    # 合肥·第十四届次元之门动漫游戏博览会 (was row 9, now row 10)
    $ws.Range("F10").Value = 6275

    # 合肥·首届Gumi同人展 (was row 12, now row 13)
    $ws.Range("F13").Value = 1265

    # 合肥·第七届环形宇宙动漫游戏嘉年华 (was row 21, now row 22)
    $ws.Range("F22").Value = 4617

    # 合肥·排球少年only之夏日招新季 (was row 22, now row 23 on "展览";
    # "全部类型" has one extra pre-existing row -> now row 24)
    if ($sheetName -eq "展览") {
        $ws.Range("F23").Value = 63
    } else {
        $ws.Range("F24").Value = 63
    }

    # 合肥·比翼连枝国乙&代号鸢only
    if ($sheetName -eq "展览") {
        $ws.Range("F24").Value = 42
    } else {
        $ws.Range("F25").Value = 42
    }

    # 合肥·第八届环形宇宙动漫游戏嘉年华Plus
    if ($sheetName -eq "展览") {
        $ws.Range("F25").Value = 74
    } else {
        $ws.Range("F26").Value = 74
    }

    # 合肥·SSS第五人格only (last row of the sheet)
    if ($sheetName -eq "展览") {
        $ws.Range("F27").Value = 74
    } else {
        $ws.Range("F28").Value = 74
    }
}
